$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 9 for the "affiliate marketing" keyword entry,
# pushing the existing rows 9-19 down to 10-20.
$ws.Rows.Item(9).Insert()
$ws.Range("A9").Value = "affiliate marketing"
$ws.Range("B9").Value = "affiliate.marketing.guide"

# Append one more "affiliate marketing" row at the very end of the table
# (new last row, 21), copying the formatting from the row above it so the
# new cells keep the same style as the rest of the table.
$ws.Range("A20:B20").Copy()
$ws.Range("A21:B21").PasteSpecial(-4122)
$ws.Range("A21").Value = "affiliate marketing"
$ws.Range("B21").Value = "affiliate.marketing.guide"

# Match the author's final selection/view state.
[void]$ws.Range("A9:B9").Select()
